$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column E ("ZoneLetter") values in rows 2-467 change from "V" to "T"
$ws.Range("E2:E467").Value = "T"
